$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database Testing")

# TC08 row (row 9): add a new Comments/Notes entry
$ws.Range("I9").Value = "This need to be checked as for Electric Vehicle user probably will not enter the EngineSize so it can be null."

# TC09 row (row 10): update Actual Result, clear Comments/Notes
$ws.Range("G10").Value = "The insertion should be successful and value is rounded of to 2 decimal places."
$ws.Range("I10").Value = ""

# TC10 row (row 11): update Expected Result
$ws.Range("F11").Value = "Record inserted successfully with Availability=1"

# Update the active selection to match the target state
$ws.Range("G18").Select()
